$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update F2, F3, F5, F6, F8, F9
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 510
$wsExhibition.Range("F3").Value = 6071
$wsExhibition.Range("F5").Value = 85
$wsExhibition.Range("F6").Value = 116
$wsExhibition.Range("F8").Value = 64
$wsExhibition.Range("F9").Value = 557

# Sheet "全部类型" (sheet4): update F2, F3, F6, F7, F10, F11
$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F2").Value = 510
$wsAllTypes.Range("F3").Value = 6071
$wsAllTypes.Range("F6").Value = 85
$wsAllTypes.Range("F7").Value = 116
$wsAllTypes.Range("F10").Value = 64
$wsAllTypes.Range("F11").Value = 557

$wb.Save()
